$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their text data type (matches source formatting)
# by explicitly setting a text number format before assigning values.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.129.75'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.71%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.053.08'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.24%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '248.12'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.89%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.665'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.64%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '57.12'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +3.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.383'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.66%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0784'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +2.96%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.96%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '16.14'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +3.29%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.912'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +17.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.349.91'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.97%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.73'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +6.20%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.066.57'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.78%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.65'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +15.85%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.207.90'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '74.86'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.88%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0892'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +2.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.45'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +3.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '238.80'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +2.54%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.50'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +6.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.95'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +3.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.56'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +6.87%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.17'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.27'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.95%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.125'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.95%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.20'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +13.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.16'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +6.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.69'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +9.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0624'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +3.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0887'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.82%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.87'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +6.16%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.31'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +7.05%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.96%  '
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'HuobiToken'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.10'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +6.28%  '
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'THORChain'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.12'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +4.70%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '102.52'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +9.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0225'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.63%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'Cronos'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0980'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -4.23%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.17'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +6.00%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '17.30'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.18%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.97%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.301.26'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +3.56%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.87'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.20%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'FTXToken'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.70'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +27.84%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.84'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +4.28%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.237.75'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.94%  '
